# Generate Report for Handback
# -----------------------------
# For both locale sheets (zh-cn, de-de) the localization round-trip has
# finished: the "Ready for handoff" status becomes "Handed back: in sync
# with en-US", the handback datetime is stamped, and the "Latest Target
# File" / "Latest Handback File" columns (I/J) are now populated with a
# link to the source doc and the name of the generated handback file.
# Column widths are widened to fit the new content.

$wb = $excel.ActiveWorkbook

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55792540f2923c0fcc3108a0d38a783322b26c2b/e2e/971a9384-d2ed-4b2e-b399-547c8512ca7b.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55792540f2923c0fcc3108a0d38a783322b26c2b/e2e/a4fc9b30-0d1e-495e-a482-ca587c3d247f.md"
$mdName1 = "971a9384-d2ed-4b2e-b399-547c8512ca7b.md"
$mdName2 = "a4fc9b30-0d1e-495e-a482-ca587c3d247f.md"

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: status text for both locales on both rows.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $handedBack
$overview.Range("F2").Value = $handedBack
$overview.Range("E3").Value = $handedBack
$overview.Range("F3").Value = $handedBack

$overview.Columns.Item(5).ColumnWidth = 29.15
$overview.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $handedBack
$zh.Range("C3").Value = $handedBack

$zh.Range("J2").Value = "971a9384-d2ed-4b2e-b399-547c8512ca7b.21b209f8499a946610a66125692f433b74ef3e68.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-15 16:44:48"

$zh.Range("J3").Value = "a4fc9b30-0d1e-495e-a482-ca587c3d247f.9b66847d11965721a4afceacf739ff219cd2849f.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-15 16:44:48"

$zh.Hyperlinks.Add($zh.Range("I2"), $mdUrl1, "", "", $mdName1)
$zh.Hyperlinks.Add($zh.Range("I3"), $mdUrl2, "", "", $mdName2)

$zh.Columns.Item(3).ColumnWidth = 29.15
$zh.Columns.Item(9).ColumnWidth = 39.15
$zh.Columns.Item(10).ColumnWidth = 39.15

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $handedBack
$de.Range("C3").Value = $handedBack

$de.Range("J2").Value = "971a9384-d2ed-4b2e-b399-547c8512ca7b.21b209f8499a946610a66125692f433b74ef3e68.de-de.xlf"
$de.Range("K2").Value = "2016-08-15 16:44:56"

$de.Range("J3").Value = "a4fc9b30-0d1e-495e-a482-ca587c3d247f.9b66847d11965721a4afceacf739ff219cd2849f.de-de.xlf"
$de.Range("K3").Value = "2016-08-15 16:44:56"

$de.Hyperlinks.Add($de.Range("I2"), $mdUrl1, "", "", $mdName1)
$de.Hyperlinks.Add($de.Range("I3"), $mdUrl2, "", "", $mdName2)

$de.Columns.Item(3).ColumnWidth = 29.15
$de.Columns.Item(9).ColumnWidth = 39.15
$de.Columns.Item(10).ColumnWidth = 39.15
